$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Spiky boi"
$ws.Range("D9").Value = "Aloe Vera"
$ws.Range("D7").Value = "Aloe Vera, Mini palm"
$ws.Range("D1").Value = "starred"

$null = $ws.Range("D1").Select()
